$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/percentage/link/name assignments ---
$ws.Cells.Item(2,4).Value = "29.389.43"
$ws.Cells.Item(2,5).Value = "  +0.09%  "
$ws.Cells.Item(3,4).Value = "1.848.60"
$ws.Cells.Item(3,5).Value = "  +0.20%  "
$ws.Cells.Item(4,5).Value = "  +0.17%  "
$ws.Cells.Item(5,5).Value = "  +0.04%  "
$ws.Cells.Item(6,5).Value = "  -0.78%  "
$ws.Cells.Item(7,5).Value = "  +0.09%  "
$ws.Cells.Item(8,5).Value = "  +0.62%  "
$ws.Cells.Item(9,5).Value = "  -1.04%  "
$ws.Cells.Item(10,5).Value = "  -1.00%  "
$ws.Cells.Item(12,4).Value = "1.852.37"
$ws.Cells.Item(12,5).Value = "  -6.68%  "
$ws.Cells.Item(13,2).Value = "Polkadot"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13,5).Value = "  +0.35%  "
$ws.Cells.Item(14,2).Value = "ShibaInu"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(14,5).Value = "  +8.97%  "
$ws.Cells.Item(15,5).Value = "  -0.72%  "
$ws.Cells.Item(16,5).Value = "  +0.59%  "
$ws.Cells.Item(17,4).Value = "2.101.82"
$ws.Cells.Item(17,5).Value = "  -7.18%  "
$ws.Cells.Item(18,5).Value = "  +0.16%  "
$ws.Cells.Item(19,4).Value = "29.400.85"
$ws.Cells.Item(19,5).Value = "  -0.02%  "
$ws.Cells.Item(20,5).Value = "  -0.89%  "
$ws.Cells.Item(21,5).Value = "  -0.10%  "
$ws.Cells.Item(22,5).Value = "  +0.10%  "
$ws.Cells.Item(23,5).Value = "  -0.83%  "
$ws.Cells.Item(24,5).Value = "  +0.18%  "
$ws.Cells.Item(25,5).Value = "  +0.63%  "
$ws.Cells.Item(26,5).Value = "  -0.73%  "
$ws.Cells.Item(27,5).Value = "  -0.26%  "
$ws.Cells.Item(28,5).Value = "  -0.32%  "
$ws.Cells.Item(29,5).Value = "  -0.23%  "
$ws.Cells.Item(30,5).Value = "  +3.92%  "
$ws.Cells.Item(31,5).Value = "  -2.10%  "
$ws.Cells.Item(32,5).Value = "  -0.44%  "
$ws.Cells.Item(33,5).Value = "  -0.02%  "
$ws.Cells.Item(34,5).Value = "  -0.60%  "
$ws.Cells.Item(35,5).Value = "  -0.01%  "
$ws.Cells.Item(36,5).Value = "  -0.84%  "
$ws.Cells.Item(37,5).Value = "  -0.39%  "
$ws.Cells.Item(38,4).Value = "1.238.31"
$ws.Cells.Item(38,5).Value = "  -0.90%  "
$ws.Cells.Item(39,5).Value = "  -0.63%  "
$ws.Cells.Item(40,5).Value = "  -1.08%  "
$ws.Cells.Item(41,5).Value = "  +5.22%  "
$ws.Cells.Item(42,5).Value = "  +0.24%  "
$ws.Cells.Item(43,5).Value = "  +0.09%  "
$ws.Cells.Item(44,5).Value = "  -0.20%  "
$ws.Cells.Item(45,5).Value = "  -0.58%  "
$ws.Cells.Item(46,5).Value = "  +2.33%  "
$ws.Cells.Item(47,5).Value = "  +1.08%  "
$ws.Cells.Item(48,5).Value = "  +0.02%  "
$ws.Cells.Item(49,5).Value = "  -1.41%  "
$ws.Cells.Item(50,5).Value = "  -1.49%  "
$ws.Cells.Item(51,5).Value = "  -0.45%  "

# --- Numeric-looking text in column D: force Text format, then reset style ---
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "1.000"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "240.15"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "0.6293"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.07592"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.2928"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "24.44"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "5.001"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.00001077"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "0.6780"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "83.60"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "6.171"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "228.34"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "7.479"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "1.002"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "157.23"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "0.1395"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "8.342"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "17.61"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "1.465"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0.05586"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "4.027"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "1.841"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.7093"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "2.587"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.01800"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "2.771"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "6.411"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.9036"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "101.64"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "65.95"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "7.144"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.4017"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "9.021"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "1.676"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "0.1119"
$ws.Cells.Item(51,4).Style = "Normal"
